$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 534-543), matching columns:
# A = WindowClassName, B = windowControlID, C = Module, D = Text
$rows = @(
    @("Edit", 22554, "Avtal", "Avtalstid från och med"),
    @("Edit", 22555, "Avtal", "Avtalstid till och med"),
    @("Edit", 22316, "Avtal", "Orderdatum"),
    @("Edit", 22556, "Avtal", "Första faktureringsmånad"),
    @("Edit", 22554, "Avtal", "Avtalstid från och med"),
    @("Edit", 22555, "Avtal", "Avtalstid till och med"),
    @("Edit", 22556, "Avtal", "Första faktureringsmånad"),
    @("Edit", 22558, "Avtal", "Fakturadag"),
    @("Edit", 22557, "Avtal", "Fakturaintervall, månader"),
    @("Edit", 22564, "Avtal", "Nästa period, slut")
)

$startRow = 534
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$ws.Activate()
$ws.Range("A543").Select()
